$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column C (rows 2-494) holds the "Förändrad" (changed) date as a serial
# number. Every row currently has the same date value 45178 and needs to
# be bumped by one day to 45179.
$ws.Range("C2:C494").Value = 45179
